# "added notes to Student B"
# Fill in the note1-note4 columns (D, F, H, J) for the studentB sheet with
# coding annotations, and update the workbook/sheet view state to reflect
# that studentB is now the active tab with the selection scrolled down to
# the bottom of the sheet where the edits were made.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("studentB")

# --- Cell value updates (note1/note2/note3/note4 columns, plus a couple of
#     theme columns that were filled in alongside them) -------------------
$values = @{
    "G3"  = "efficiency"
    "H3"  = "reusing code instead of rewriting it"
    "D4"  = "descriptive comment on process taken in code below"
    "D5"  = "full path to data"
    "D6"  = "full path to data"
    "D7"  = "filtering observations based on inclusion operator (%in%) and logical comparison (!)"
    "D8"  = "mutate variable"
    "D9"  = "descriptive comment on actions taken in code below"
    "D11" = "select variable"
    "G11" = "data structures"
    "H11" = "vector"
    "E12" = "data wrangling"
    "F12" = "variable selection"
    "D13" = "descriptive comment on actions taken in code below"
    "D14" = "variable selection"
    "D15" = "code comment on units of calculation above"
    "D16" = "descriptive comment on actions taken in code below"
    "D17" = "create function to estimate parameter"
    "E17" = "data structures"
    "F17" = "vector"
    "D18" = "code comment marking new section of code"
    "D19" = "descriptive comment on actions taken in code below"
    "D20" = "select variable, filter rows"
    "D21" = "select variable, filter rows"
    "D22" = "filtering observations based on logical comparison"
    "D23" = "filtering observations based on logical comparison"
    "D24" = "descriptive comment on action taken in code below"
    "D25" = "filter observations"
    "D26" = "descriptive comment on action taken in code below"
    "D27" = "select variable, filter rows"
    "D28" = "descriptive comment on action taken in code below"
    "D30" = "descriptive comment on action taken in code below"
    "H32" = "pull out MLE estimate"
    "D35" = "code comment on units of calculation below"
    "H36" = "obtain point estimate"
    "D37" = "descriptive comment on action taken in code below"
    "A39" = "par(mar = c(3.5, 4, 3, 1))"
    "F41" = "obtain predictions"
    "E46" = "data structures"
    "F46" = "vector"
    "D47" = "descriptive comment on process in code below"
    "D48" = "descriptive comment on action taken in code below"
    "D52" = "matrix"
    "D53" = "applying function across columns of matrix"
    "E53" = "efficiency"
    "F53" = "repeated operations on multiple rows"
    "H53" = "obtaining likelihood estimates"
    "F54" = "obtaining minimum of likelihood estimates"
    "F55" = "locating which index corresponds to minimum"
    "D57" = "filtering observations based on index"
    "F57" = "matrix"
    "J57" = "obtaining lowerbound for confidence interval"
    "D58" = "filtering observations based on index"
    "F58" = "matrix"
    "J58" = "obtaining upperbound for confidence interval"
    "F59" = "obtaining confidence interval"
    "D60" = "inspect object"
    "F61" = "adjusting confidence interval lowerbound"
    "F62" = "adjusting confidence interval upperbound"
    "D63" = "descriptive comment on process taken in code below"
    "E65" = "data structures"
    "F65" = "vector"
    "C67" = "data visualization"
}

foreach ($addr in $values.Keys) {
    $ws2.Range($addr).Value = $values[$addr]
}

# Cells that had notes before but are now cleared out
$ws2.Range("D66").Value = ""
$ws2.Range("J59").Value = ""

# A handful of the newly-typed note cells came in with an explicit black
# font color rather than the sheet's usual theme color -- match that.
$blackFontCells = @("F4", "D16", "D19", "D37", "D47", "D48", "D63")
foreach ($addr in $blackFontCells) {
    $ws2.Range($addr).Font.Color = 0
}

# --- View / selection state ------------------------------------------------
# studentA was the active tab before, scrolled to A38 with A45 selected.
# After these edits studentB becomes the active tab, scrolled near the
# bottom of the sheet (row 37) with A66 selected.
$ws1 = $wb.Worksheets.Item("studentA")
$ws1.Select()
$ws1.Range("A45").Select()

$ws2.Select()
$ws2.Range("A66").Select()
